$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (Förändrad) increments from 46064 to 46065 for all data rows 2-16
for ($r = 2; $r -le 16; $r++) {
    $ws.Cells.Item($r, 3).Value = 46065
}

# Rows 7-16: reorder A/B/G values according to the new arrangement
$ws.Cells.Item(7, 1).Value = "A 50530-2024"
$ws.Cells.Item(7, 2).Value = 45601.56424768519
$ws.Cells.Item(7, 7).Value = 0.7

$ws.Cells.Item(8, 1).Value = "A 50538-2024"
$ws.Cells.Item(8, 2).Value = 45601.57153935185
$ws.Cells.Item(8, 7).Value = 0.8

$ws.Cells.Item(9, 1).Value = "A 23677-2023"
$ws.Cells.Item(9, 2).Value = 45077
$ws.Cells.Item(9, 7).Value = 0.6

$ws.Cells.Item(10, 1).Value = "A 2253-2022"
$ws.Cells.Item(10, 2).Value = 44578
$ws.Cells.Item(10, 7).Value = 0.3

$ws.Cells.Item(11, 1).Value = "A 11351-2021"
$ws.Cells.Item(11, 2).Value = 44263
$ws.Cells.Item(11, 7).Value = 0.5

$ws.Cells.Item(12, 1).Value = "A 32633-2025"
$ws.Cells.Item(12, 2).Value = 45838.65677083333
$ws.Cells.Item(12, 7).Value = 1.3

$ws.Cells.Item(13, 1).Value = "A 50277-2024"
$ws.Cells.Item(13, 2).Value = 45600.60440972223
$ws.Cells.Item(13, 7).Value = 0.5

$ws.Cells.Item(14, 1).Value = "A 23678-2023"
$ws.Cells.Item(14, 2).Value = 45077
$ws.Cells.Item(14, 7).Value = 1.4

$ws.Cells.Item(15, 1).Value = "A 45370-2022"
$ws.Cells.Item(15, 2).Value = 44844.6397337963
$ws.Cells.Item(15, 7).Value = 2.7

$ws.Cells.Item(16, 1).Value = "A 58926-2025"
$ws.Cells.Item(16, 2).Value = 45986
$ws.Cells.Item(16, 7).Value = 3.1
